# LOM3236.xlsx update
# - Bumps the "Ativação:" date from 01/01/2016 to 01/01/2023 (shown in
#   both the B/C "Ativação:" row and the B/C "Programa resumido:" row,
#   which happen to share the same text in the source workbook)
# - Fills in the English "Objectives", "Short syllabus" and "Syllabus"
#   descriptions (columns B/C) that were previously blank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-TextValue($Cell, $Text, $FormatSourceCell) {
    # Force the cell to Text so a value that looks like a date (e.g.
    # "01/01/2023") is kept as literal text instead of being converted to
    # a date serial number.
    $ws.Range($Cell).NumberFormat = "@"
    $ws.Range($Cell).Value = $Text
    # Re-apply the normal column formatting (font/alignment/General number
    # format) from a neighbouring cell that already has it, so the cell
    # ends up styled the same way as the rest of the column instead of
    # being stuck with the temporary "@" text format.
    $ws.Range($FormatSourceCell).Copy() | Out-Null
    $ws.Range($Cell).PasteSpecial($xlPasteFormats) | Out-Null
}

# 1) Update activation date everywhere it appears (B8/C8 "Ativação:" and
#    B13/C13 "Programa resumido:", which mirror the same text)
Set-TextValue "B8"  "01/01/2023" "B9"
Set-TextValue "C8"  "01/01/2023" "C9"
Set-TextValue "B13" "01/01/2023" "B9"
Set-TextValue "C13" "01/01/2023" "C9"

# 2) Objectives (row 11) - English objectives text
Set-TextValue "B11" "Develop theoretical and practical knowledge of the manufacturing processes of equipment and devices required for the development of products and prototypes. Know the requirements and effects of manufacturing processes in order to allow, interact, create and execute projects throughout your professional life." "B10"
Set-TextValue "C11" "Develop theoretical and practical knowledge of the manufacturing processes of equipment and devices required for the development of products and prototypes. Know the requirements and effects of manufacturing processes in order to allow, interact, create and execute projects throughout your professional life." "C10"

# 3) Short syllabus (row 14) - English short syllabus text
Set-TextValue "B14" "Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping." "B10"
Set-TextValue "C14" "Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping." "C10"

# 4) Syllabus (row 16) - English syllabus text
Set-TextValue "B16" "Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder)." "B10"
Set-TextValue "C16" "Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder)." "C10"

$excel.CutCopyMode = 0
